$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing data (D:K) to (F:M)
$ws.Range("D:E").Insert()

# Copy number formats/styles from the (now-shifted) old columns into the new D:E
# columns so the new cells inherit the same per-row formatting (date format on
# the "Period Ending" rows, number format elsewhere). Limited to the row blocks
# that actually contain data so blank separator rows stay untouched.
$ws.Range("F7:M35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:M77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:M102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns with the new quarter's data.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 923900
$ws.Range("E8").Value = 970300
$ws.Range("D9").Value = 671000
$ws.Range("E9").Value = 701800
$ws.Range("D10").Value = 252900
$ws.Range("E10").Value = 268500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 802900
$ws.Range("E17").Value = 833500
$ws.Range("D18").Value = 121000
$ws.Range("E18").Value = 136800
$ws.Range("D20").Value = -37800
$ws.Range("E20").Value = -40900
$ws.Range("D21").Value = 188900
$ws.Range("E21").Value = 198600
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 83100
$ws.Range("E23").Value = 95900
$ws.Range("D24").Value = 22900
$ws.Range("E24").Value = 25000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 60300
$ws.Range("E26").Value = 70900
$ws.Range("D27").Value = 59100
$ws.Range("E27").Value = 69300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 25300
$ws.Range("E29").Value = 3700
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 37800
$ws.Range("E32").Value = 40900
$ws.Range("D33").Value = 84400
$ws.Range("E33").Value = 73000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 84400
$ws.Range("E35").Value = 73000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 509300
$ws.Range("E41").Value = 444300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 598300
$ws.Range("E43").Value = 589900
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 80000
$ws.Range("E45").Value = 170100
$ws.Range("D46").Value = 1187600
$ws.Range("E46").Value = 1204300
$ws.Range("D47").Value = 27800
$ws.Range("E47").Value = 28700
$ws.Range("D48").Value = 790400
$ws.Range("E48").Value = 788000
$ws.Range("D49").Value = 3165600
$ws.Range("E49").Value = 3182700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 635000
$ws.Range("E52").Value = 664800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 5806400
$ws.Range("E54").Value = 5868500
$ws.Range("D57").Value = 165200
$ws.Range("E57").Value = 158800
$ws.Range("D58").Value = 68400
$ws.Range("E58").Value = 64200
$ws.Range("D59").Value = 784700
$ws.Range("E59").Value = 829200
$ws.Range("D60").Value = 1018400
$ws.Range("E60").Value = 1052200
$ws.Range("D61").Value = 3337500
$ws.Range("E61").Value = 3355600
$ws.Range("D62").Value = 476200
$ws.Range("E62").Value = 522800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 4839300
$ws.Range("E66").Value = 4936800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -768600
$ws.Range("E72").Value = -814400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 967100
$ws.Range("E76").Value = 931700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 84400
$ws.Range("E81").Value = 73000
$ws.Range("D83").Value = 105800
$ws.Range("E83").Value = 102700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 186100
$ws.Range("E89").Value = 198100
$ws.Range("D91").Value = -78300
$ws.Range("E91").Value = -73800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -69600
$ws.Range("E94").Value = -73800
$ws.Range("D96").Value = -38500
$ws.Range("E96").Value = -38500
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -54100
$ws.Range("E100").Value = -50900
$ws.Range("D101").Value = 2600
$ws.Range("E101").Value = 900
$ws.Range("D102").Value = 64900
$ws.Range("E102").Value = 74300

$ws.Calculate()
